# Update report header text (Volume/Number and week-covering dates)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Volume 31   Number  20"
$ws.Range("C9").Value = "Report Covering the Week  5/13/2024  Through  5/19/2024"

# Update weekly crime-statistics grid (rows 14-30)
$ws.Range("G14").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("G14").PasteSpecial(-4122)
$ws.Range("H14").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("H14").PasteSpecial(-4122)
$ws.Range("C15").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("E15").Value = -100
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -33.333333333333
$ws.Range("I15").Value = 4
$ws.Range("J15").Value = 7
$ws.Range("K15").Value = -42.857142857142
$ws.Range("L15").Value = -55.555555555555
$ws.Range("M15").Value = -20
$ws.Range("N15").Value = -85.714285714285
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -40
$ws.Range("F16").Value = 18
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = 12.5
$ws.Range("I16").Value = 80
$ws.Range("J16").Value = 67
$ws.Range("K16").Value = 19.402985074626
$ws.Range("L16").Value = 25
$ws.Range("M16").Value = -11.111111111111
$ws.Range("N16").Value = -72.972972972973
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = 14.285714285714
$ws.Range("F17").Value = 36
$ws.Range("G17").Value = 34
$ws.Range("H17").Value = 5.882352941176
$ws.Range("I17").Value = 162
$ws.Range("J17").Value = 134
$ws.Range("K17").Value = 20.895522388059
$ws.Range("L17").Value = 35
$ws.Range("M17").Value = 131.428571428571
$ws.Range("N17").Value = -22.857142857142
$ws.Range("D15").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -80
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -27.272727272727
$ws.Range("I18").Value = 68
$ws.Range("J18").Value = 46
$ws.Range("K18").Value = 47.826086956521
$ws.Range("L18").Value = 6.25
$ws.Range("M18").Value = 195.652173913043
$ws.Range("N18").Value = -55.555555555555
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 12.5
$ws.Range("F19").Value = 34
$ws.Range("G19").Value = 25
$ws.Range("H19").Value = 36
$ws.Range("I19").Value = 145
$ws.Range("J19").Value = 158
$ws.Range("K19").Value = -8.227848101265
$ws.Range("L19").Value = -5.844155844155
$ws.Range("M19").Value = 72.619047619047
$ws.Range("N19").Value = -36.403508771929
$ws.Range("D15").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 33.333333333333
$ws.Range("I20").Value = 22
$ws.Range("K20").Value = -40.540540540540
$ws.Range("L20").Value = -42.105263157894
$ws.Range("M20").Value = 29.411764705882
$ws.Range("N20").Value = -84.397163120567
$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = -11.538461538461
$ws.Range("F21").Value = 102
$ws.Range("G21").Value = 92
$ws.Range("H21").Value = 10.869565217391
$ws.Range("I21").Value = 483
$ws.Range("J21").Value = 455
$ws.Range("K21").Value = 6.153846153846
$ws.Range("L21").Value = 6.858407079646
$ws.Range("M21").Value = 66.551724137931
$ws.Range("N21").Value = -54.901960784313
$ws.Range("G22").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("H22").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("H22").PasteSpecial(-4122)
$ws.Range("C23").Value = 7
$ws.Range("D23").Value = 8
$ws.Range("E23").Value = -12.5
$ws.Range("I23").Value = 138
$ws.Range("J23").Value = 139
$ws.Range("K23").Value = -0.719424460431
$ws.Range("L23").Value = -4.827586206896
$ws.Range("M23").Value = 68.292682926829
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = 0
$ws.Range("I24").Value = 303
$ws.Range("J24").Value = 344
$ws.Range("K24").Value = -11.918604651162
$ws.Range("L24").Value = 8.602150537634
$ws.Range("M24").Value = 29.487179487179
$ws.Range("C25").Value = 6
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 21
$ws.Range("G25").Value = 29
$ws.Range("H25").Value = -27.586206896551
$ws.Range("I25").Value = 64
$ws.Range("J25").Value = 87
$ws.Range("K25").Value = -26.436781609195
$ws.Range("L25").Value = 4.918032786885
$ws.Range("C26").Value = 17
$ws.Range("D26").Value = 16
$ws.Range("E26").Value = 6.25
$ws.Range("F26").Value = 66
$ws.Range("G26").Value = 44
$ws.Range("H26").Value = 50
$ws.Range("I26").Value = 275
$ws.Range("J26").Value = 208
$ws.Range("K26").Value = 32.211538461538
$ws.Range("L26").Value = 27.314814814814
$ws.Range("M26").Value = 8.267716535433
$ws.Range("C27").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -50
$ws.Range("I27").Value = 8
$ws.Range("J27").Value = 13
$ws.Range("K27").Value = -38.461538461538
$ws.Range("L27").Value = -27.272727272727
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("G28").PasteSpecial(-4122)
$ws.Range("H28").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("H28").PasteSpecial(-4122)
$ws.Range("D29").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("G29").Value = 2
$ws.Range("D30").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("G30").Value = 1

$excel.CutCopyMode = $false
